$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The historical series gained two more (older) year-end data points at the
# top, pushing the previously existing rows down by two. Shift the existing
# data rows (old rows 2-8, now rows 4-10) down manually, from the bottom up,
# so that copied values don't get clobbered before they're read, and so we
# avoid the extra formatting side effects that Rows.Insert() would bring in.

for ($r = 8; $r -ge 2; $r--) {
    $destRow = $r + 2
    $ws.Cells.Item($destRow, 1).Value = $ws.Cells.Item($r, 1).Value2
    $ws.Cells.Item($destRow, 2).Value = $ws.Cells.Item($r, 2).Value2
}

# New row 2: 2014-12-31 -> 921462000000
$ws.Cells.Item(2, 1).Value = 42004
$ws.Cells.Item(2, 2).Value = 921462000000

# New row 3: 2015-12-31 -> 820805000000
$ws.Cells.Item(3, 1).Value = 42369
$ws.Cells.Item(3, 2).Value = 820805000000

# Apply the same formatting used by the rest of the date column (the custom
# "YYYY-MM-DD HH:MM:SS" number format with centered/top alignment and thin
# border) to the two newly introduced A-column cells, as well as the two
# cells at the bottom of the column that are now populated for the first
# time and would otherwise have no explicit style.
$ws.Range("A4").Copy()
$ws.Range("A2:A3").PasteSpecial(-4122) | Out-Null
$ws.Range("A9:A10").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = $false
